$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 5.81
$ws.Range("H2").Value = 7.17
$ws.Range("I2").Value = 5.81
$ws.Range("J2").Value = 5.95
$ws.Range("K2").Value = 2.32
$ws.Range("L2").Value = 7.17
$ws.Range("M2").Value = 4.29

# Row 3
$ws.Range("E3").Value = 6.07
$ws.Range("H3").Value = 6.96
$ws.Range("I3").Value = 6.07
$ws.Range("J3").Value = 5.78
$ws.Range("K3").Value = 2.61
$ws.Range("L3").Value = 6.96
$ws.Range("M3").Value = 4.46

# Row 4
$ws.Range("E4").Value = 5.91
$ws.Range("H4").Value = 7.05
$ws.Range("I4").Value = 5.91
$ws.Range("J4").Value = 5.85
$ws.Range("K4").Value = 2.42
$ws.Range("L4").Value = 7.05
$ws.Range("M4").Value = 4.35

# Row 5
$ws.Range("E5").Value = 5.7
$ws.Range("H5").Value = 7.26
$ws.Range("I5").Value = 5.7
$ws.Range("J5").Value = 6.02
$ws.Range("K5").Value = 2.36
$ws.Range("L5").Value = 7.26
$ws.Range("M5").Value = 4.32

# Row 6
$ws.Range("E6").Value = 5.76
$ws.Range("H6").Value = 7.15
$ws.Range("I6").Value = 5.76
$ws.Range("J6").Value = 5.93
$ws.Range("K6").Value = 2.23
$ws.Range("L6").Value = 7.15
$ws.Range("M6").Value = 4.24

# Row 7
$ws.Range("E7").Value = 5.64
$ws.Range("I7").Value = 5.64
$ws.Range("J7").Value = 6.17
$ws.Range("K7").Value = 2.26
$ws.Range("M7").Value = 4.26
